$d = $word.ActiveDocument

# --- Simple in-place text replacements ---

# Customer name changed; new text keeps a trailing space -> xml:space="preserve"
$d.Content.Find.Execute("Customer: Tigo Bolivia", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Customer:JuanchoMix ", 2) | Out-Null

# Date changed
$d.Content.Find.Execute("Date & Time: 09/07/20", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Date & Time: 09/22/20", 2) | Out-Null

# Cliente ID changed (bold run, text only)
$d.Content.Find.Execute("Cliente ID: 232323", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Cliente ID: 111", 2) | Out-Null

# Root Cause filled in with a value
$d.Content.Find.Execute("Root Cause: ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Root Cause: QMI.json", 2) | Out-Null

# --- Append a new block of paragraphs at the end of the document ---
# (repeats the "Cliente ID / Descripcion / ---- / 2. DETAILS ... / Root Cause" block
#  for a second customer id, as plain (non-bold) paragraphs first so that formatting
#  doesn't leak between paragraphs via InsertParagraphAfter)

$newTexts = @(
    "Cliente ID: 222",
    "Descripcion: ",
    "--------------------------------------------------------------------",
    "2. DETAILS OF INCIDENT: ",
    "Impacted platform: ",
    "Root Cause: LatePacket.json",
    "Incident description: ",
    "Evidencias: ",
    "--------------------------------------------------------------------",
    "3. RESOLUTION",
    "Incident Analysis: ",
    "Workaround: ",
    "Recommendation: ",
    "Additional comments: NA"
)
$boldFlags = @($true, $true, $false, $false, $false, $false, `
               $false, $false, $false, $false, $false, $false, $false, $false)

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$lastPara.Range.InsertParagraphAfter()

$idx = $count + 1
for ($i = 0; $i -lt $newTexts.Count; $i++) {
    $text = $newTexts[$i]
    $p = $d.Paragraphs($idx)
    $p.Range.Text = $text
    if ($idx -lt ($count + $newTexts.Count)) {
        $p2 = $d.Paragraphs($idx)
        $p2.Range.InsertParagraphAfter()
    }
    $idx = $idx + 1
}

# Apply bold formatting only to the runs that need it (not the paragraph mark)
$idx = $count + 1
for ($i = 0; $i -lt $newTexts.Count; $i++) {
    if ($boldFlags[$i]) {
        $text = $newTexts[$i]
        $p = $d.Paragraphs($idx)
        $r = $p.Range
        $boldRange = $d.Range($r.Start, $r.Start + $text.Length)
        $boldRange.Font.Bold = $true
    }
    $idx = $idx + 1
}
